$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price data between row 2 and row 3
# (Date, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)

$ws.Range("D2").Value = 44839
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15600
$ws.Range("P2").Value = 1040

$ws.Range("D3").Value = 44750
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 19000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 19571
$ws.Range("P3").Value = 1305
